$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 11.83761344122534
$ws.Range("C2").Value = 11.00155019027585
$ws.Range("E2").Value = 12.71584953030981
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 29.19859257109157
$ws.Range("H2").Value = 14.59936334985623
$ws.Range("K2").Value = 8.240296873008354
$ws.Range("L2").Value = 9.485956752826516
$ws.Range("O2").Value = 22.24187474788997

# Row 3
$ws.Range("B3").Value = 11.52843795131523
$ws.Range("C3").Value = 11.03809021072026
$ws.Range("E3").Value = 12.73696528279488
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 29.39869849743449
$ws.Range("H3").Value = 14.66093121901194
$ws.Range("K3").Value = 8.003531196468225
$ws.Range("L3").Value = 9.467616026833815
$ws.Range("O3").Value = 22.35962859529359

# Row 4
$ws.Range("B4").Value = 11.33557391989793
$ws.Range("C4").Value = 11.06192649742061
$ws.Range("E4").Value = 12.75252592020799
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 29.53061050858843
$ws.Range("H4").Value = 14.70095236703374
$ws.Range("K4").Value = 7.854890680898942
$ws.Range("L4").Value = 9.457872896246771
$ws.Range("O4").Value = 22.43648184753093

# Row 5
$ws.Range("B5").Value = 11.25633047948722
$ws.Range("C5").Value = 11.07199292503824
$ws.Range("E5").Value = 12.75951934589662
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 29.58663417217993
$ws.Range("H5").Value = 14.71781992696755
$ws.Range("K5").Value = 7.793579404771825
$ws.Range("L5").Value = 9.454287314174604
$ws.Range("O5").Value = 22.46894506644487

# Row 6
$ws.Range("B6").Value = 11.24313631267511
$ws.Range("C6").Value = 11.07368578739623
$ws.Range("E6").Value = 12.76071998628646
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 29.5960736866469
$ws.Range("H6").Value = 14.72065453610075
$ws.Range("K6").Value = 7.783356589185848
$ws.Range("L6").Value = 9.453715257688451
$ws.Range("O6").Value = 22.47440471725459

# Row 7
$ws.Range("B7").Value = 11.3345076936793
$ws.Range("C7").Value = 11.0620608266907
$ws.Range("E7").Value = 12.75261759522182
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 29.53135688664745
$ws.Range("H7").Value = 14.70117758591992
$ws.Range("K7").Value = 7.854066695294525
$ws.Range("L7").Value = 9.45782297784457
$ws.Range("O7").Value = 22.43691502181175

# Row 8
$ws.Range("B8").Value = 11.73170357652461
$ws.Range("C8").Value = 11.01385896024576
$ws.Range("E8").Value = 12.72259140595024
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 29.26570775334213
$ws.Range("H8").Value = 14.62013205455302
$ws.Range("K8").Value = 8.159385219062415
$ws.Range("L8").Value = 9.479319701434131
$ws.Range("O8").Value = 22.28153168635237

# Row 9
$ws.Range("B9").Value = 12.48198393289533
$ws.Range("C9").Value = 10.93041375836401
$ws.Range("E9").Value = 12.68431264016031
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 28.81683289606589
$ws.Range("H9").Value = 14.47876490638842
$ws.Range("K9").Value = 8.728900668508876
$ws.Range("L9").Value = 9.53337454553011
$ws.Range("O9").Value = 22.01293010550816

# Row 10
$ws.Range("B10").Value = 13.01011617072371
$ws.Range("C10").Value = 10.87581293186725
$ws.Range("E10").Value = 12.66875238506102
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 28.53137236058527
$ws.Range("H10").Value = 14.38555639562637
$ws.Range("K10").Value = 9.125516645975924
$ws.Range("L10").Value = 9.580131585533318
$ws.Range("O10").Value = 21.83758125749166

# Row 11
$ws.Range("B11").Value = 13.24430988101386
$ws.Range("C11").Value = 10.85241990326206
$ws.Range("E11").Value = 12.6643987217668
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 28.41122295159762
$ws.Range("H11").Value = 14.34545589136057
$ws.Range("K11").Value = 9.300486643299859
$ws.Range("L11").Value = 9.60288135460741
$ws.Range("O11").Value = 21.76258350015423

# Row 12
$ws.Range("B12").Value = 13.33204411693814
$ws.Range("C12").Value = 10.84376864495668
$ws.Range("E12").Value = 12.66314133456098
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 28.36712929619947
$ws.Range("H12").Value = 14.33060092921311
$ws.Range("K12").Value = 9.365906360711765
$ws.Range("L12").Value = 9.611704183123445
$ws.Range("O12").Value = 21.73486970369989

# Row 13
$ws.Range("B13").Value = 13.31319240594157
$ws.Range("C13").Value = 10.84562264415381
$ws.Range("E13").Value = 12.66339474558472
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 28.37656306215178
$ws.Range("H13").Value = 14.33378553313477
$ws.Range("K13").Value = 9.35185510236769
$ws.Range("L13").Value = 9.609794858028389
$ws.Range("O13").Value = 21.74080783834748

# Row 14
$ws.Range("B14").Value = 13.2515472927307
$ws.Range("C14").Value = 10.85170401057665
$ws.Range("E14").Value = 12.6642874384809
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 28.40756714715667
$ws.Range("H14").Value = 14.34422714966984
$ws.Range("K14").Value = 9.305885837390001
$ws.Range("L14").Value = 9.603603074071525
$ws.Range("O14").Value = 21.76028971387694

# Row 15
$ws.Range("B15").Value = 13.21366196128395
$ws.Range("C15").Value = 10.85545598451401
$ws.Range("E15").Value = 12.66488517077698
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 28.42674117788169
$ws.Range("H15").Value = 14.35066592988958
$ws.Range("K15").Value = 9.277617755963332
$ws.Range("L15").Value = 9.59983736665497
$ws.Range("O15").Value = 21.77231230530243

# Row 16
$ws.Range("B16").Value = 12.99468228181471
$ws.Range("C16").Value = 10.87737074679892
$ws.Range("E16").Value = 12.66909169252977
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 28.53942046423421
$ws.Range("H16").Value = 14.38822329979568
$ws.Range("K16").Value = 9.113967522221216
$ws.Range("L16").Value = 9.578674206210406
$ws.Range("O16").Value = 21.84257854118074

# Row 17
$ws.Range("B17").Value = 12.85873574508749
$ws.Range("C17").Value = 10.89118442616292
$ws.Range("E17").Value = 12.67236979554604
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 28.61103787652776
$ws.Range("H17").Value = 14.41185232496258
$ws.Range("K17").Value = 9.012137011228411
$ws.Range("L17").Value = 9.566066955752383
$ws.Range("O17").Value = 21.88690640341778

# Row 18
$ws.Range("B18").Value = 12.7799771082631
$ws.Range("C18").Value = 10.89926575204245
$ws.Range("E18").Value = 12.67451178887952
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 28.65314357444472
$ws.Range("H18").Value = 14.42565968401112
$ws.Range("K18").Value = 8.953056412091641
$ws.Range("L18").Value = 9.558955183888838
$ws.Range("O18").Value = 21.91285149921022

# Row 19
$ws.Range("B19").Value = 12.75321613518515
$ws.Range("C19").Value = 10.90202533870006
$ws.Range("E19").Value = 12.67528110183901
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 28.66755644240859
$ws.Range("H19").Value = 14.43037183256578
$ws.Range("K19").Value = 8.932966799486968
$ws.Range("L19").Value = 9.556571376295166
$ws.Range("O19").Value = 21.92171314751062

# Row 20
$ws.Range("B20").Value = 12.87326663696576
$ws.Range("C20").Value = 10.88969985907961
$ws.Range("E20").Value = 12.67199429102287
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 28.60331950239379
$ws.Range("H20").Value = 14.40931456513365
$ws.Range("K20").Value = 9.023030270887629
$ws.Range("L20").Value = 9.567394606557064
$ws.Range("O20").Value = 21.88214116649357

# Row 21
$ws.Range("B21").Value = 13.26968032934643
$ws.Range("C21").Value = 10.84991214823552
$ws.Range("E21").Value = 12.66401462002445
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 28.39842231591956
$ws.Range("H21").Value = 14.34115123600843
$ws.Range("K21").Value = 9.319411257413675
$ws.Range("L21").Value = 9.605416145949263
$ws.Range("O21").Value = 21.75454878722939

# Row 22
$ws.Range("B22").Value = 13.52319151040626
$ws.Range("C22").Value = 10.82511585919242
$ws.Range("E22").Value = 12.66107954738689
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 28.27269958036476
$ws.Range("H22").Value = 14.29852713044328
$ws.Range("K22").Value = 9.5209080838242
$ws.Range("L22").Value = 9.631475316647281
$ws.Range("O22").Value = 21.67515988139119

# Row 23
$ws.Range("B23").Value = 13.3884215718642
$ws.Range("C23").Value = 10.83823984710665
$ws.Range("E23").Value = 12.6624376604391
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 28.33904804601441
$ws.Range("H23").Value = 14.32110050604629
$ws.Range("K23").Value = 9.407909511865622
$ws.Range("L23").Value = 9.617457990196968
$ws.Range("O23").Value = 21.71716507460704

# Row 24
$ws.Range("B24").Value = 12.86669908716123
$ws.Range("C24").Value = 10.89037059679143
$ws.Range("E24").Value = 12.67216325491056
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 28.60680607718358
$ws.Range("H24").Value = 14.4104611925537
$ws.Range("K24").Value = 9.018107096341993
$ws.Range("L24").Value = 9.566793950698878
$ws.Range("O24").Value = 21.88429409596327

# Row 25
$ws.Range("B25").Value = 12.28268847572323
$ws.Range("C25").Value = 10.9518068076248
$ws.Range("E25").Value = 12.6924613654821
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 28.93051036271637
$ws.Range("H25").Value = 14.51513378470868
$ws.Range("K25").Value = 8.578395344688493
$ws.Range("L25").Value = 9.517498072188058
$ws.Range("O25").Value = 22.08173076231023
